# Corrects vr calculation (per Rio): welfare is now computed against
# local GDP (column F, gdp_pc_pp) instead of national GDP (column K,
# gdp_pc_pp_nat). This updates the 'v_r' (asset vulnerability, non-poor
# people) values in column M for rows 4-88.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    4 = 0.15242128150753;
    5 = 0.14724254636513;
    6 = 0.16585666659158;
    7 = 0.206958140009178;
    8 = 0.151701628531602;
    9 = 0.181891571459366;
    10 = 0.125986370157819;
    11 = 0.130239329179526;
    12 = 0.171087536883518;
    13 = 0.107909775226284;
    14 = 0.216726125137212;
    15 = 0.117115029621362;
    16 = 0.102282370153621;
    17 = 0.157837245475113;
    18 = 0.153647306539452;
    19 = 0.126189427700763;
    20 = 0.116596190575278;
    21 = 0.130822279753631;
    22 = 0.166789710102733;
    23 = 0.163813413930527;
    24 = 0.123894990624163;
    25 = 0.198294078972812;
    26 = 0.16071094980118;
    27 = 0.108798222756118;
    28 = 0.153195889727578;
    29 = 0.148883898312743;
    30 = 0.151198378718521;
    31 = 0.14080956109018;
    32 = 0.163991603714574;
    33 = 0.133144632755999;
    34 = 0.171221791961544;
    35 = 0.1542390449639;
    36 = 0.176569582409104;
    37 = 0.11599043098252;
    38 = 0.119258567384796;
    39 = 0.113781081637163;
    40 = 0.17152392811848;
    41 = 0.134264699279154;
    42 = 0.156583373378184;
    43 = 0.127818871208354;
    44 = 0.121431696379744;
    45 = 0.112658594228783;
    46 = 0.121874471324649;
    47 = 0.130257361938168;
    48 = 0.153015407242403;
    49 = 0.220262356843528;
    50 = 0.112270576613113;
    51 = 0.144329833534141;
    52 = 0.193198501872659;
    53 = 0.17835053914082;
    54 = 0.154132021167951;
    55 = 0.105328963340617;
    56 = 0.112404011864711;
    57 = 0.113205689192865;
    58 = 0.108188043797672;
    59 = 0.17187298757124;
    60 = 0.179002925547794;
    61 = 0.204158284530675;
    62 = 0.126151476726511;
    63 = 0.11108931349105;
    64 = 0.188862504264756;
    65 = 0.155836855895197;
    66 = 0.249450371632919;
    67 = 0.115122358895172;
    68 = 0.137066537445382;
    69 = 0.160280920637414;
    70 = 0.109007161915461;
    71 = 0.109811425134469;
    72 = 0.178461879432624;
    73 = 0.180139223905763;
    74 = 0.199160191128915;
    75 = 0.144174295565331;
    76 = 0.18881789715651;
    77 = 0.15545908567038;
    78 = 0.137361819211783;
    79 = 0.175768657603379;
    80 = 0.26537811714166;
    81 = 0.185302806499261;
    82 = 0.180356271635341;
    83 = 0.124954676418872;
    84 = 0.15930999050001;
    85 = 0.126452227767051;
    86 = 0.165919371829714;
    87 = 0.190818921327535;
    88 = 0.190569579519451
}

foreach ($row in $newValues.Keys) {
    $ws.Range("M$row").Value = $newValues[$row]
}
